$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings: _old -> _FV2410, _new -> _FV2504
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2410")
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2504")
}

# Freeze the header row (pane split after row 1)
$selected = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a proper Excel Table ("Table1") with an autofilter
$dataRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
